$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6002.1
$ws.Range("J17").Value = 6926
$ws.Range("L17").Value = 20778
$ws.Range("N17").Value = -21114
$ws.Range("H33").Value = 297.27274
$ws.Range("I33").Value = 274.77777
$ws.Range("J33").Value = 398.5
$ws.Range("K33").Value = 274.77777
$ws.Range("L33").Value = 398.5
$ws.Range("M33").Value = -45.77776999999998
$ws.Range("N33").Value = -856.5
$ws.Range("H34").Value = 9066.666999999999
$ws.Range("I34").Value = 5966.6665
$ws.Range("J34").Value = 12166.667
$ws.Range("K34").Value = 5966.6665
$ws.Range("L34").Value = 12166.667
$ws.Range("M34").Value = -5763.6665
$ws.Range("N34").Value = -12572.667
$ws.Range("H36").Value = 9066.666999999999
$ws.Range("I36").Value = 5966.6665
$ws.Range("J36").Value = 12166.667
$ws.Range("K36").Value = 5966.6665
$ws.Range("L36").Value = 12166.667
$ws.Range("M36").Value = -5251.6665
$ws.Range("N36").Value = -13596.667
$ws.Range("H41").Value = 2864.9
$ws.Range("I41").Value = 2849.4285
$ws.Range("J41").Value = 2901
$ws.Range("K41").Value = 2849.4285
$ws.Range("L41").Value = 2901
$ws.Range("M41").Value = -2409.4285
$ws.Range("N41").Value = -3781
$ws.Range("H62").Value = 6023.696
$ws.Range("I62").Value = 5275.077
$ws.Range("J62").Value = 6996.9
$ws.Range("K62").Value = 5275.077
$ws.Range("L62").Value = 6996.9
$ws.Range("M62").Value = -4651.077
$ws.Range("N62").Value = -8244.9
$ws.Range("H65").Value = 6023.696
$ws.Range("I65").Value = 5275.077
$ws.Range("J65").Value = 6996.9
$ws.Range("K65").Value = 26375.385
$ws.Range("L65").Value = 34984.5
$ws.Range("M65").Value = -23255.385
$ws.Range("N65").Value = -41224.5
$ws.Range("H69").Value = 7905.591
$ws.Range("I69").Value = 2004.3334
$ws.Range("J69").Value = 8837.368
$ws.Range("K69").Value = 6013.0002
$ws.Range("L69").Value = 26512.104
$ws.Range("M69").Value = -5139.0002
$ws.Range("N69").Value = -28260.104
$ws.Range("H72").Value = 7905.591
$ws.Range("I72").Value = 2004.3334
$ws.Range("J72").Value = 8837.368
$ws.Range("K72").Value = 18039.0006
$ws.Range("L72").Value = 79536.31200000001
$ws.Range("M72").Value = -13671.0006
$ws.Range("N72").Value = -88272.31200000001
$ws.Range("H96").Value = 1300.9375
$ws.Range("I96").Value = 1405.3572
$ws.Range("K96").Value = 4216.071599999999
$ws.Range("M96").Value = -2843.071599999999
$ws.Range("H111").Value = 1760
$ws.Range("I111").Value = 833.3333
$ws.Range("K111").Value = 2499.9999
$ws.Range("M111").Value = 567.0001000000002
$ws.Range("H113").Value = 7141.7856
$ws.Range("J113").Value = 7707.1665
$ws.Range("L113").Value = 7707.1665
$ws.Range("N113").Value = -14215.1665
$ws.Range("H137").Value = 2349.4119
$ws.Range("I137").Value = 1389.3846
$ws.Range("J137").Value = 2943.7144
$ws.Range("K137").Value = 4168.1538
$ws.Range("L137").Value = 8831.143199999999
$ws.Range("M137").Value = -1618.1538
$ws.Range("N137").Value = -13931.1432
$ws.Range("H141").Value = 1055.3334
$ws.Range("J141").Value = 1000
$ws.Range("L141").Value = 3000
$ws.Range("N141").Value = -13360
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1765.902
$ws.Range("I32").Value = 1233.0212
$ws.Range("K32").Value = 1233.0212
$ws.Range("M32").Value = -946.0211999999999
$ws.Range("H45").Value = 5496.6
$ws.Range("I45").Value = 1842
$ws.Range("K45").Value = 1842
$ws.Range("M45").Value = -1465
$ws.Range("H74").Value = 15153024
$ws.Range("I74").Value = 15874501
$ws.Range("K74").Value = 15874501
$ws.Range("M74").Value = -15873627
$ws.Range("H77").Value = 15153024
$ws.Range("I77").Value = 15874501
$ws.Range("K77").Value = 79372505
$ws.Range("M77").Value = -79368137
$ws.Range("H122").Value = 2545.3333
$ws.Range("I122").Value = 1614.4445
$ws.Range("K122").Value = 4843.333500000001
$ws.Range("M122").Value = -2393.333500000001
$ws.Range("H131").Value = 54195
$ws.Range("J131").Value = 54195
$ws.Range("L131").Value = 54195
$ws.Range("N131").Value = -64275
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3643.3845
$ws.Range("I20").Value = 2838.1333
$ws.Range("J20").Value = 4741.4546
$ws.Range("K20").Value = 2838.1333
$ws.Range("L20").Value = 4741.4546
$ws.Range("M20").Value = -2591.1333
$ws.Range("N20").Value = -5235.4546
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4108.077
$ws.Range("J16").Value = 5242.2
$ws.Range("L16").Value = 5242.2
$ws.Range("N16").Value = -5816.2
$ws.Range("H31").Value = 33751.8
$ws.Range("I31").Value = 4395.2593
$ws.Range("K31").Value = 4395.2593
$ws.Range("M31").Value = -4100.2593
$ws.Range("H34").Value = 33751.8
$ws.Range("I34").Value = 4395.2593
$ws.Range("K34").Value = 4395.2593
$ws.Range("M34").Value = -4193.2593
$ws.Range("H113").Value = 4108.077
$ws.Range("J113").Value = 5242.2
$ws.Range("L113").Value = 5242.2
$ws.Range("N113").Value = -9582.200000000001
$ws.Range("H122").Value = 6574.0557
$ws.Range("I122").Value = 2568
$ws.Range("J122").Value = 16989.8
$ws.Range("K122").Value = 7704
$ws.Range("L122").Value = 50969.39999999999
$ws.Range("M122").Value = -5254
$ws.Range("N122").Value = -55869.39999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 11000.75
$ws.Range("J42").Value = 11000.75
$ws.Range("L42").Value = 33002.25
$ws.Range("N42").Value = -34070.25
$ws.Range("H56").Value = 6075.5557
$ws.Range("I56").Value = 6075.5557
$ws.Range("K56").Value = 6075.5557
$ws.Range("M56").Value = -5545.5557
$ws.Range("H112").Value = 83337950
$ws.Range("I112").Value = 166674000
$ws.Range("K112").Value = 500022000
$ws.Range("M112").Value = -500020892
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9879.777
$ws.Range("I70").Value = 6652.5
$ws.Range("J70").Value = 16334.333
$ws.Range("K70").Value = 6652.5
$ws.Range("L70").Value = 16334.333
$ws.Range("M70").Value = -6382.5
$ws.Range("N70").Value = -16874.333
$ws.Range("H73").Value = 9879.777
$ws.Range("I73").Value = 6652.5
$ws.Range("J73").Value = 16334.333
$ws.Range("K73").Value = 6652.5
$ws.Range("L73").Value = 16334.333
$ws.Range("M73").Value = -5716.5
$ws.Range("N73").Value = -18206.333
$ws.Range("H101").Value = 22200
$ws.Range("J101").Value = 22200
$ws.Range("L101").Value = 22200
$ws.Range("N101").Value = -28690
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 15899.8
$ws.Range("J101").Value = 15899.8
$ws.Range("L101").Value = 15899.8
$ws.Range("N101").Value = -22389.8
$ws.Range("H130").Value = 65954.89
$ws.Range("J130").Value = 65954.89
$ws.Range("L130").Value = 65954.89
$ws.Range("N130").Value = -75994.89
$ws.Range("H132").Value = 4942.4116
$ws.Range("I132").Value = 2967.5833
$ws.Range("J132").Value = 9682
$ws.Range("K132").Value = 8902.749899999999
$ws.Range("L132").Value = 29046
$ws.Range("M132").Value = -6372.749899999999
$ws.Range("N132").Value = -34106
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4556.2144
$ws.Range("I81").Value = 2898.75
$ws.Range("K81").Value = 5797.5
$ws.Range("M81").Value = -4736.5
$ws.Range("H84").Value = 4556.2144
$ws.Range("I84").Value = 2898.75
$ws.Range("K84").Value = 28987.5
$ws.Range("M84").Value = -23683.5
$ws.Range("H123").Value = 47625
$ws.Range("J123").Value = 47625
$ws.Range("L123").Value = 47625
$ws.Range("N123").Value = -57425
$ws.Range("H132").Value = 11203.714
$ws.Range("I132").Value = 8643.454
$ws.Range("J132").Value = 14020
$ws.Range("K132").Value = 25930.362
$ws.Range("L132").Value = 42060
$ws.Range("M132").Value = -23400.362
$ws.Range("N132").Value = -47120
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
